$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("C1").Value = "Big StockRoom"

# New "Cooler" column header, matching the existing header style (bold/centered/bordered)
$ws.Range("C1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Cooler"
$ws.Application.CutCopyMode = $false

# --- Row 2: Hashbrowns ---
$ws.Range("A2").Value = "Hashbrowns"
$ws.Range("B2").Value = 50
$ws.Range("C2").ClearContents()
$ws.Range("E2").Value = "Y"

# --- Row 3: Oil ---
$ws.Range("A3").Value = "Oil"
$ws.Range("B3").Value = 10

# --- Row 4: Ketchup ---
$ws.Range("A4").Value = "Ketchup"
$ws.Range("B4").Value = 5

# --- Remove old rows 5-7 (data now only spans rows 1-4) ---
$ws.Range("A5:E7").Delete()
